$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91:232 down to 92:233
$ws.Rows("91").Insert()

# Populate the newly inserted row 91 with the new data point
$ws.Range("A91").Value = 5
$ws.Range("B91").Value = "Macroferia Regional de Talca"
$ws.Range("C91").Value = "Maule"
$ws.Range("D91").Value = 44557
$ws.Range("E91").Value = 7
$ws.Range("F91").Value = 100114014
$ws.Range("G91").Value = "Betarraga"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 550
$ws.Range("L91").Value = 550
$ws.Range("M91").Value = 550
$ws.Range("N91").Value = "`$/paquete 5 unidades"
$ws.Range("O91").Value = "Región del Maule"
$ws.Range("P91").Value = 110
$ws.Range("Q91").Value = 5
$ws.Range("R91").Value = "Hortaliza"
